$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the header in column E (shared string "Ampel/Kreuzung?" -> "Ampel?")
$ws.Range("E1").Value = "Ampel?"

# Mark "DWPT-Abschnitt?" (column D) as active (1) for rows 3 through 10
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1

# Update the speed value in B6
$ws.Range("B6").Value = 20

# Move the active selection to G9 (matches recorded view state)
$ws.Range("G9").Select()
